{"js": "// Change the \"Source Code\" console-output styling so sessionInfo() results\n// fit horizontally: shrink the font size from 9.5pt to 9pt on both the\n// \"Source Code\" paragraph style and its linked \"Verbatim Char\" character\n// style.\n\nconst styles = context.document.getStyles();\n\nconst verbatimChar = styles.getByName(\"Verbatim Char\");\nverbatimChar.font.load(\"size\");\n\nconst sourceCode = styles.getByName(\"Source Code\");\nsourceCode.font.load(\"size\");\n\nawait context.sync();\n\nverbatimChar.font.size = 9;\nsourceCode.font.size = 9;\n\nawait context.sync();\n", "ps1": "# Change the \"Source Code\" console-output styling so sessionInfo() results\n# fit horizontally: shrink the font size from 9.5pt to 9pt on both the\n# \"Source Code\" paragraph style and its linked \"Verbatim Char\" character\n# style.\n\n$d = $word.ActiveDocument\n\n$verbatimChar = $d.Styles.Item(\"Verbatim Char\")\n$verbatimChar.Font.Size = 9\n\n$sourceCode = $d.Styles.Item(\"Source Code\")\n$sourceCode.Font.Size = 9\n"}
